# Change R10 to 430k to increase minimum release time.
# R10 used to share the 100k / Axial 6.8mm / MCMF006FF1003A50 BoM line with
# R43, R29 and R1. It now gets its own BoM line using a 530k (430k marking)
# metal-film resistor, and is dropped from the shared 100k line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for R10 right above the former R13/R4 (56R) row, which
# shifts every subsequent row (and the shared-formula / SUM ranges) down by
# one — Excel handles that automatically.
$ws.Rows.Item(22).Insert()

# Pick up the number formats / alignment of the row it was inserted above
# (left-aligned value column, text order-code column, currency price/total
# columns) so the new row matches the rest of the table instead of using
# bare defaults. (Only A:K — this BoM line has no Notes/L entry.)
$ws.Range("A23:K23").Copy()
$ws.Range("A22:K22").PasteSpecial(-4122)

$ws.Range("A22").Value = "R10"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "430k"
$ws.Range("D22").Value = "Axial 6.8mm"
$ws.Range("E22").Value = "530k 0.6W 1% metal film"
$ws.Range("F22").Value = "Multicomp"
$ws.Range("G22").Value = "MCMF006FF4303A50"
$ws.Range("H22").Value = "Farnell"
$ws.Range("I22").Value = "2401823"
$ws.Range("J22").Value = 0.0207
$ws.Range("K22").Formula = "=B22*J22"

# The old combined "R43 R29 R10 R1" / 100k row is now row 34 (it was row 33
# before the insert above). Drop R10 from the reference list and its
# quantity, since it now has its own line.
$ws.Range("A34").Value = "R43 R29 R1 "
$ws.Range("B34").Value = 3
